# Auto-generated edit script applying cached-value corrections
# to the Asura_Profits leve-profit columns (H:N) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3803.2222
$ws.Range("I63").Value = 5501.6665
$ws.Range("J63").Value = 2954
$ws.Range("K63").Value = 5501.6665
$ws.Range("L63").Value = 2954
$ws.Range("M63").Value = -4815.6665
$ws.Range("N63").Value = -4326

$ws.Range("H66").Value = 3803.2222
$ws.Range("I66").Value = 5501.6665
$ws.Range("J66").Value = 2954
$ws.Range("K66").Value = 27508.3325
$ws.Range("L66").Value = 14770
$ws.Range("M66").Value = -24076.3325
$ws.Range("N66").Value = -21634

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 129686.25
$ws.Range("I20").Value = 334730
$ws.Range("J20").Value = 6660
$ws.Range("K20").Value = 334730
$ws.Range("L20").Value = 6660
$ws.Range("M20").Value = -334483
$ws.Range("N20").Value = -7154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 25002296
$ws.Range("I22").Value = 100000400
$ws.Range("J22").Value = 2929.1667
$ws.Range("K22").Value = 300001200
$ws.Range("L22").Value = 8787.500100000001
$ws.Range("M22").Value = -300001031
$ws.Range("N22").Value = -9125.500100000001

$ws.Range("H27").Value = 25002296
$ws.Range("I27").Value = 100000400
$ws.Range("J27").Value = 2929.1667
$ws.Range("K27").Value = 300001200
$ws.Range("L27").Value = 8787.500100000001
$ws.Range("M27").Value = -300001098
$ws.Range("N27").Value = -8991.500100000001

$ws.Range("H35").Value = 9860
$ws.Range("J35").Value = 12250
$ws.Range("L35").Value = 36750
$ws.Range("N35").Value = -37326

$ws.Range("H41").Value = 845.65216
$ws.Range("J41").Value = 870.4545000000001
$ws.Range("L41").Value = 2611.3635
$ws.Range("N41").Value = -3287.3635

$ws.Range("H64").Value = 4531.773
$ws.Range("I64").Value = 1740
$ws.Range("J64").Value = 5352.8823
$ws.Range("K64").Value = 5220
$ws.Range("L64").Value = 16058.6469
$ws.Range("M64").Value = -4950
$ws.Range("N64").Value = -16598.6469

$ws.Range("H67").Value = 4531.773
$ws.Range("I67").Value = 1740
$ws.Range("J67").Value = 5352.8823
$ws.Range("K67").Value = 5220
$ws.Range("L67").Value = 16058.6469
$ws.Range("M67").Value = -4284
$ws.Range("N67").Value = -17930.6469

$ws.Range("H113").Value = 590.73914
$ws.Range("I113").Value = 525.1667
$ws.Range("J113").Value = 826.8
$ws.Range("K113").Value = 1575.5001
$ws.Range("L113").Value = 2480.4
$ws.Range("M113").Value = 594.4999
$ws.Range("N113").Value = -6820.4

$ws.Range("H141").Value = 5732
$ws.Range("I141").Value = 6464
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 19392
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -14212
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1219.5217
$ws.Range("I126").Value = 790.5
$ws.Range("J126").Value = 1687.5454
$ws.Range("K126").Value = 2371.5
$ws.Range("L126").Value = 5062.6362
$ws.Range("M126").Value = 98.5
$ws.Range("N126").Value = -10002.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 506.03705
$ws.Range("I107").Value = 503.5263
$ws.Range("J107").Value = 512
$ws.Range("K107").Value = 1510.5789
$ws.Range("L107").Value = 1536
$ws.Range("M107").Value = 409.4211
$ws.Range("N107").Value = -5376

$ws.Range("H119:N119").ClearContents()

$ws.Range("H120:N120").ClearContents()

$ws.Range("H121:N121").ClearContents()

$ws.Range("H122:N122").ClearContents()

$ws.Range("H123:N123").ClearContents()

$ws.Range("H124:N124").ClearContents()

$ws.Range("H125:N125").ClearContents()

$ws.Range("H126:N126").ClearContents()

$ws.Range("H127:N127").ClearContents()

$ws.Range("H128:N128").ClearContents()

$ws.Range("H129:N129").ClearContents()

$ws.Range("H130:N130").ClearContents()

$ws.Range("H131:N131").ClearContents()

$ws.Range("H132:N132").ClearContents()

$ws.Range("H133:N133").ClearContents()

$ws.Range("H135:N135").ClearContents()

$ws.Range("H136:N136").ClearContents()

$ws.Range("H137:N137").ClearContents()

$ws.Range("H138:N138").ClearContents()

$ws.Range("H139:N139").ClearContents()

$ws.Range("H140:N140").ClearContents()

$ws.Range("H141:N141").ClearContents()
